# Updates the cryptos list (columns D = Price, E = Volume(1h)) with fresh
# values, matching the commit "Updated cryptos list ... with GitHub Actions".
#
# Each cell in D/E is stored in the workbook as a literal (inline) text
# string - e.g. "27.734.68" or "  -0.29%  " - not a real number, so when we
# write the new value back we force the cell format to Text first and then
# restore the Normal style afterwards (this avoids Excel auto-converting a
# value like "1.012" into the number 1.012, while also avoiding a stray
# "quote prefix" style being left behind on the cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = "27.734.68"; E = "  -0.29%  " }
    @{ Row = 3; D = "1.850.48"; E = "  -0.90%  " }
    @{ Row = 4; D = "1.012"; E = "  -2.63%  " }
    @{ Row = 5; D = "319.66"; E = "  -1.56%  " }
    @{ Row = 6; D = "1.012"; E = "  -2.39%  " }
    @{ Row = 7; D = "0.4338"; E = "  -2.07%  " }
    @{ Row = 8; D = $null; E = "  -0.78%  " }
    @{ Row = 9; D = "0.07390"; E = "  -1.15%  " }
    @{ Row = 10; D = "0.8841"; E = "  -0.33%  " }
    @{ Row = 11; D = "21.67"; E = "  -0.65%  " }
    @{ Row = 12; D = "1.859.23"; E = "  -1.08%  " }
    @{ Row = 13; D = "6.754"; E = "  -0.19%  " }
    @{ Row = 14; D = "5.476"; E = "  -1.58%  " }
    @{ Row = 15; D = "0.07141"; E = "  -1.39%  " }
    @{ Row = 16; D = "88.22"; E = "  +5.12%  " }
    @{ Row = 17; D = $null; E = "  -2.48%  " }
    @{ Row = 18; D = "0.000009041"; E = "  -1.43%  " }
    @{ Row = 19; D = "1.012"; E = "  -2.42%  " }
    @{ Row = 20; D = "15.55"; E = "  -0.20%  " }
    @{ Row = 21; D = "27.746.99"; E = "  -0.39%  " }
    @{ Row = 22; D = "5.266"; E = "  -1.16%  " }
    @{ Row = 23; D = "11.23"; E = "  -1.41%  " }
    @{ Row = 24; D = "2.083.73"; E = "  -1.73%  " }
    @{ Row = 25; D = "2.029"; E = "  +2.82%  " }
    @{ Row = 26; D = "155.75"; E = "  -1.94%  " }
    @{ Row = 27; D = "18.65"; E = "  -1.37%  " }
    @{ Row = 28; D = "2.149"; E = "  +7.72%  " }
    @{ Row = 29; D = "5.423"; E = "  +1.64%  " }
    @{ Row = 30; D = "120.57"; E = "  +2.32%  " }
    @{ Row = 31; D = "0.08969"; E = "  -1.55%  " }
    @{ Row = 32; D = "1.238"; E = "  +1.34%  " }
    @{ Row = 33; D = "0.7777"; E = "  -0.40%  " }
    @{ Row = 34; D = "4.582"; E = "  -0.15%  " }
    @{ Row = 35; D = "2.922"; E = "  -5.59%  " }
    @{ Row = 36; D = "1.143"; E = "  -1.91%  " }
    @{ Row = 37; D = "1.013"; E = "  -2.66%  " }
    @{ Row = 38; D = "0.05348"; E = "  -0.33%  " }
    @{ Row = 39; D = "0.01975"; E = "  -1.40%  " }
    @{ Row = 40; D = "7.193"; E = "  +4.22%  " }
    @{ Row = 41; D = "2.865"; E = "  +0.27%  " }
    @{ Row = 42; D = "0.5192"; E = "  -0.45%  " }
    @{ Row = 43; D = "0.1686"; E = "  -0.69%  " }
    @{ Row = 44; D = $null; E = "  +2.63%  " }
    @{ Row = 45; D = "110.93"; E = "  +0.79%  " }
    @{ Row = 46; D = "10.78"; E = "  +0.48%  " }
    @{ Row = 48; D = "0.4751"; E = "  +0.58%  " }
    @{ Row = 49; D = "0.06507"; E = "  +0.75%  " }
    @{ Row = 50; D = "1.013"; E = "  -2.57%  " }
    @{ Row = 51; D = "1.905"; E = "  +0.08%  " }
)

foreach ($u in $updates) {
    $row = $u.Row

    if ($null -ne $u.D) {
        $dCell = $ws.Cells.Item($row, 4)
        $dCell.NumberFormat = "@"
        $dCell.Value = $u.D
        $dCell.Style = "Normal"
    }

    $eCell = $ws.Cells.Item($row, 5)
    $eCell.NumberFormat = "@"
    $eCell.Value = $u.E
    $eCell.Style = "Normal"
}
